# Update automatico via Actualizar 09-09-2020 00-41-10
# Adds a new row (row 41) to the "Trabajo_CL32" table on sheet "trabajo"
# describing the second-debate approval of the food-voucher law amendment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trabajo")

$tbl = $ws.ListObjects.Item("Trabajo_CL32")

# Grow the table by one row (table ref A1:K40 -> A1:K41)
$newListRow = $tbl.ListRows.Add()

# Copy the formatting (styles) of the last existing data row (40) into the
# freshly added row (41) so fonts/fills/borders/number formats match.
$ws.Range("A40:K40").Copy()
$ws.Range("A41:K41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height matches the other wrapped-text rows (120pt, same as row 40).
$ws.Rows.Item(41).RowHeight = 120

# --- Cell values for the new row -------------------------------------------
$ws.Range("A41").Value = "Ministerio de Trabajo y Desarrollo Social"
$ws.Range("B41").Value = 40
$ws.Range("C41").Value = "Trabajo"
$ws.Range("D41").Value = "El Ministerio de Trabajo y Desarrollo Laboral de Panamá (MITRADEL) es un Ministerio de la República de Panamá que forma parte del Órgano Ejecutivo."
$ws.Range("E41").Value = "https://www.mitradel.gob.pa/modificaciones-a-la-ley-de-vale-alimenticio-aprobadas-en-segundo-debate/"
$ws.Range("F41").Value = "Durante la sesión del martes 8 de septiembre de 2020, el Pleno de la Asamblea Nacional de Diputados aprobó en segundo debate el proyecto de ley 365, mediante el cual se modifica la Ley 59 del 7 de agosto de 2003, sobre el Programa de Alimentación de Trabajadores."
$ws.Range("G41").Value = "https://www.mitradel.gob.pa"
$ws.Range("H41").Value = 44082
$ws.Range("I41").Value = 44082
$ws.Range("J41").Value = "Panamá"
$ws.Range("K41").Value = "Ministerial"

# --- Hyperlinks for the new row ---------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G41"), "https://www.mitradel.gob.pa/")
$ws.Hyperlinks.Add($ws.Range("E41"), "https://www.mitradel.gob.pa/modificaciones-a-la-ley-de-vale-alimenticio-aprobadas-en-segundo-debate/")

# --- Data validation range (Categoria dropdown) extended to new row --------
$dv = $ws.Range("C2:C40").Validation
$dv.Delete()
$ws.Range("C2:C41").Validation.Add(0)
$dv2 = $ws.Range("C2:C41").Validation
$dv2.ErrorTitle = "Entrada no válida"
$dv2.ErrorMessage = "Selecciona una categoría de la lista"
$dv2.InputTitle = "Categoria"
$dv2.InputMessage = "Selecciona una categoría de la lista"
$dv2.ShowInput = $true
$dv2.ShowError = $true

# --- Selection / view position moved to the newly added row ----------------
$ws.Activate()
$ws.Range("F41").Select()

Write-Host "Row 41 added to table Trabajo_CL32"
